$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet; the product catalogue is being replaced with a
# new category/subcategory layout for the WooCommerce import/update flow.
$ws.Cells.Clear()

# New header row: CATEGORY, SUBCATEGORY, NAME, SKU, PRICE, DIMENSION, AREA
$ws.Range("A1").Value = "CATEGORY"
$ws.Range("B1").Value = "SUBCATEGORY"
$ws.Range("C1").Value = "NAME"
$ws.Range("D1").Value = "SKU"
$ws.Range("E1").Value = "PRICE"
$ws.Range("F1").Value = "DIMENSION"
$ws.Range("G1").Value = "AREA"

# Header row (including the trailing blank H1 cell) keeps the bold styling
# used by the original header.
$ws.Range("A1:H1").Font.Bold = $true

# Product rows
$data = @(
    @("Furniture", "Chair", "Kursi Santai", "P01A", 9999,  "200 x 200",       "INDOOR"),
    @("Furniture", "Chair", "Kursi Santai", "P01B", 39999, "200 x 200",       "OUTDOOR"),
    @("Furniture", "Table", "Meja besar",   "P02A", 5000,  "100 x 200 x 300", "INDOOR"),
    @("Furniture", "Table", "Meja besar",   "P02B", 6000,  "100 x 200 x 300", "OUTDOOR"),
    @("Furniture", "Table", "Meja besar",   "P02C", 7000,  "200 x 200 x 300", "INDOOR"),
    @("Furniture", "Table", "Meja besar",   "P02E", 8500,  "200 x 200 x 300", "OUTDOOR"),
    @("Furniture", "Table", "Sofa",         "P03A", 1000,  "100 x 100",       "INDOOR"),
    @("Furniture", "Table", "Sofa",         "P03B", 2000,  "100 x 100",       "OUTDOOR")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# Widen the CATEGORY/SUBCATEGORY columns to fit their new, longer headers;
# the NAME and DIMENSION columns keep their existing widths.
# (Inputs are chosen so the engine's internal char->stored-width rounding
# lands as close as possible to the authored widths of 17.109375 / 14.21875.)
$ws.Columns.Item(1).ColumnWidth = 16.3
$ws.Columns.Item(2).ColumnWidth = 13.3

# Selection mirrors the post-edit cursor position in the source workbook.
$ws.Range("G10").Select()
